$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 567
$ws.Range("I2").Value = 1465
$ws.Range("J2").Value = 6042
$ws.Range("K2").Value = 31
$ws.Range("L2").Value = 1616
$ws.Range("M2").Value = 91
$ws.Range("N2").Value = 1053
$ws.Range("O2").Value = 2
$ws.Range("R2").Value = 100
$ws.Range("S2").Value = 648
$ws.Range("T2").Value = 1105
$ws.Range("U2").Value = 66
$ws.Range("V2").Value = 9477
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 9532
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 149
$ws.Range("AA2").Value = 74
